$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 162, shifting all existing rows (162-266) down to (163-267).
$ws.Rows.Item(162).Insert()

# Populate the newly inserted row 162 with the new record's data.
$ws.Range("A162").Value = 4
$ws.Range("B162").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C162").Value = "Los Lagos"
$ws.Range("D162").Value = 44767
$ws.Range("E162").Value = 10
$ws.Range("F162").Value = 100112032
$ws.Range("G162").Value = "Zapallo italiano"
$ws.Range("H162").Value = "Sin especificar"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 70
$ws.Range("K162").Value = 17000
$ws.Range("L162").Value = 17000
$ws.Range("M162").Value = 17000
$ws.Range("N162").Value = "$/caja 50 unidades"
$ws.Range("O162").Value = "Región de Arica y Parinacota"
$ws.Range("P162").Value = 340
$ws.Range("Q162").Value = 50
$ws.Range("R162").Value = "Hortaliza"
